$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    "datastruct",
    "algo",
    "sysprog",
    "ver",
    "build",
    "test",
    "probdec",
    "sysdec",
    "com",
    "orgfile",
    "ordxfile",
    "tree",
    "read",
    "def",
    "err",
    "ide",
    "api",
    "fw",
    "req",
    "scr",
    "db",
    "lang",
    "plat",
    "yrs",
    "dom",
    "tool",
    "langexp",
    "cbaseexp",
    "upcom",
    "platint",
    "book",
    "blog"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("B15").Select()
